# Update "想去人数" (want-to-go count) figures for gh-pages output regeneration.
# Sheet 1 = "展览" (Exhibition), Sheet 4 = "全部类型" (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 763
$ws1.Range("F3").Value  = 51
$ws1.Range("F6").Value  = 115
$ws1.Range("F8").Value  = 112
$ws1.Range("F10").Value = 425
$ws1.Range("F11").Value = 496
$ws1.Range("F12").Value = 131
$ws1.Range("F13").Value = 11328
$ws1.Range("F14").Value = 5360

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 763
$ws4.Range("F3").Value  = 51
$ws4.Range("F8").Value  = 115
$ws4.Range("F10").Value = 112
$ws4.Range("F12").Value = 425
$ws4.Range("F13").Value = 496
$ws4.Range("F14").Value = 131
$ws4.Range("F15").Value = 11328
$ws4.Range("F17").Value = 5360
